$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style/format of an existing header cell (H1) onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$iVals = @(8,8,8,8,8,7,8,7,6,8,8,8,7,8,9,7,7,8,8,8,7,8,8,8,9,10,8,7,8,8,8,7,8,8,8,8,7,7,6,8,7,8,6,8,7,8,7,6,8,8,8,7,6,7,8,7,8,8,7,7,7,8,8,7,8,8,8,7,9,6,7,6,4,4,6)
$jVals = @(8,8,8,9,8,7,8,7,7,8,8,8,8,8,9,8,8,8,8,8,8,8,8,9,9,10,8,8,8,8,8,7,8,8,8,8,7,8,7,8,8,8,7,8,7,8,8,7,8,8,8,7,7,7,8,8,8,8,7,8,7,8,8,7,8,8,8,7,9,6,7,6,4,4,6)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
